$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 784
$ws.Range("C8").Value = 597.0999999999999
$ws.Range("C12").Value = 408.8
$ws.Range("C14").Value = 2631.6
$ws.Range("C15").Value = 198.6
$ws.Range("C18").Value = 29.85
$ws.Range("C20").Value = 652.3
$ws.Range("C24").Value = 829.9
$ws.Range("C27").Value = 160.4
$ws.Range("C29").Value = 9.600000000000001
$ws.Range("C30").Value = 517.5
$ws.Range("C31").Value = 743.1
$ws.Range("C37").Value = 2631.3
$ws.Range("C38").Value = 6.4
$ws.Range("C41").Value = 99.34999999999999
$ws.Range("C45").Value = 9.499999999999998
$ws.Range("C48").Value = 675.4999999999999
$ws.Range("C51").Value = 261
$ws.Range("C59").Value = 181.3
$ws.Range("C60").Value = 276.9
$ws.Range("C63").Value = 117.25
$ws.Range("C67").Value = 325.8
$ws.Range("C68").Value = 161.8
$ws.Range("C71").Value = 297.5
$ws.Range("C73").Value = 168.8
$ws.Range("C79").Value = 72.2
$ws.Range("C80").Value = 88.10000000000001
$ws.Range("C81").Value = 217.5
$ws.Range("C82").Value = 2735.900000000001
$ws.Range("C83").Value = 98.29999999999998
$ws.Range("C86").Value = 183.2
$ws.Range("C90").Value = 427.4999999999999
$ws.Range("C105").Value = 138.6
